# Add a new "2022-Q3" quarter sheet (right after the "总计" summary sheet)
# and record it in the summary table, pushing the older quarters down.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

$summary = $wb.Worksheets.Item(1)   # "总计"

# ---------------------------------------------------------------------
# 1) Shift the existing summary rows (2022-Q2 .. 2021-Q1) down by one row
#    and insert the new 2022-Q3 row at the top of the data (row 2).
# ---------------------------------------------------------------------
$existing = @()
for ($r = 2; $r -le 7; $r++) {
  $existing += ,@(
    $summary.Cells.Item($r, 2).Value2,
    $summary.Cells.Item($r, 3).Value2,
    $summary.Cells.Item($r, 4).Value2
  )
}

for ($i = 0; $i -lt $existing.Length; $i++) {
  $destRow = $i + 3
  $summary.Cells.Item($destRow, 2).Value = $existing[$i][0]
  $summary.Cells.Item($destRow, 3).Value = $existing[$i][1]
  $summary.Cells.Item($destRow, 4).Value = $existing[$i][2]
}

# Row 8 is brand new - give column A the same look as the rest of the
# index column (style copied from an existing indexed cell).
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial($xlPasteFormats)
$summary.Cells.Item(8, 1).Value = 6

# Write the new first data row: 2022-Q3
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.31

# ---------------------------------------------------------------------
# 2) Insert a brand-new worksheet "2022-Q3" right after "总计" holding
#    the per-fund holdings detail for the quarter.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
  $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "'014887"
$q3.Cells.Item(2, 3).Value = "招商安福1年定期开放债券"
$q3.Cells.Item(2, 4).Value = "'17.28"
$q3.Cells.Item(2, 5).Value = "'29.21"
$q3.Cells.Item(2, 6).Value = "'1.78"
$q3.Cells.Item(2, 7).Value = "'0.3076"
$q3.Cells.Item(2, 8).Value = 4

# The values above that look like numbers get an implicit "text" number
# format from the leading apostrophe trick - strip that back off so the
# cells end up with the default (no explicit) style, same as the sibling
# quarter sheets.
$q3.Range("B2:G2").ClearFormats()

# Re-apply the bold/centered/bordered header style used throughout the
# workbook to the header row and to the first (index) column.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats)

$summary.Range("A2").Copy()
$q3.Range("A2").PasteSpecial($xlPasteFormats)

# Restore the original active sheet ("总计") so the workbook view state
# matches what it was before the edit.
$summary.Activate()
$excel.CutCopyMode = $false
